# Update COMP2 spray drift workbook: insert a "method" column into the
# spray_drift sheet, recording whether each card was Fine / Spot / Coarse.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("spray_drift")

# Insert a new column before column C (old C:G -> new D:H), shifting the
# existing trt_orig/rep/all/yellow/percent_sprayed columns right by one.
$ws.Columns("C").Insert()
$ws.Columns("C").ColumnWidth = $ws.Columns("B").ColumnWidth

# Header for the new column.
$ws.Cells.Item(1, 3).Value = "method"

# Row -> method lookup (based on the treatment naming: FineBoom*, Spot*,
# coarseboom*).
$fineRows   = @(2,3,4,5,6,17,18,19,20,21)
$spotRows   = @(7,8,9,10,11,22,23,24,25,26)
$coarseRows = @(12,13,14,15,16,27,28,29,30,31)

# Touch the three distinct labels in the same order they were authored
# (Spot, Fine, Coarse) before filling the rest, so new shared-string
# entries land in that order.
$ws.Cells.Item(7, 3).Value = "Spot"
$ws.Cells.Item(2, 3).Value = "Fine"
$ws.Cells.Item(12, 3).Value = "Coarse"

foreach ($r in $fineRows)   { $ws.Cells.Item($r, 3).Value = "Fine" }
foreach ($r in $spotRows)   { $ws.Cells.Item($r, 3).Value = "Spot" }
foreach ($r in $coarseRows) { $ws.Cells.Item($r, 3).Value = "Coarse" }

# Match the saved view state from the edit: scrolled down with C27:C31
# selected.
$ws.Application.Goto($ws.Range("A11"), $false)
$ws.Range("C27:C31").Select()
